$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 11) down to the two new rows (12-13)
$ws.Range("A11:H11").Copy()
$ws.Range("A12:H12").PasteSpecial(-4122)
$ws.Range("A11:H11").Copy()
$ws.Range("A13:H13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 12 values
$ws.Range("A12").Value = "x"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = "db--"

# Row 13 values
$ws.Range("A13").Value = "x"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = "db++"

$ws.Range("G12").Select()
